# Apply the "Updated symbol list" refresh of cryptocurrency prices/volumes.
# All cells on this sheet are stored as text (inline strings), including
# values that look numeric/percentages, so every assignment below is
# prefixed with a leading apostrophe to force Excel to keep it as text
# (this also preserves formatting such as trailing zeros, e.g. "0.3500").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.11"
$ws.Range("E2").Value = "'-0.02%"
$ws.Range("D3").Value = "'41.66"
$ws.Range("E3").Value = "'1.31%"
$ws.Range("D4").Value = "'5.693"
$ws.Range("E4").Value = "'-0.10%"
$ws.Range("D5").Value = "'0.08424"
$ws.Range("E5").Value = "'4.45%"
$ws.Range("D6").Value = "'8.791"
$ws.Range("E7").Value = "'-1.70%"
$ws.Range("D8").Value = "'4.488"
$ws.Range("E8").Value = "'-0.54%"
$ws.Range("E9").Value = "'0.52%"
$ws.Range("D10").Value = "'0.9270"
$ws.Range("E10").Value = "'0.51%"
$ws.Range("D11").Value = "'0.1276"
$ws.Range("E11").Value = "'0.63%"
$ws.Range("D12").Value = "'0.1977"
$ws.Range("E12").Value = "'1.71%"
$ws.Range("D13").Value = "'0.09410"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("D14").Value = "'0.03925"
$ws.Range("E14").Value = "'6.28%"
$ws.Range("E15").Value = "'0.69%"
$ws.Range("D16").Value = "'0.001310"
$ws.Range("E16").Value = "'0.88%"
$ws.Range("D17").Value = "'0.006115"
$ws.Range("E17").Value = "'-2.49%"
$ws.Range("B18").Value = "'LEO"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.422"
$ws.Range("E18").Value = "'1.71%"
$ws.Range("B19").Value = "'BitpandaEcosystemToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3500"
$ws.Range("E19").Value = "'0.73%"
$ws.Range("B20").Value = "'MCDex"
$ws.Range("C20").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.960"
$ws.Range("E20").Value = "'8.30%"
$ws.Range("B21").Value = "'ProBitToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("E21").Value = "'-3.75%"
$ws.Range("B22").Value = "'ZBToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2513"
$ws.Range("E22").Value = "'-5.26%"
$ws.Range("B23").Value = "'CoinExToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("E23").Value = "'-0.47%"
$ws.Range("B24").Value = "'BitKan"
$ws.Range("C24").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("E24").Value = "'-1.18%"
$ws.Range("B25").Value = "'HotbitToken"
$ws.Range("C25").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004378"
$ws.Range("E25").Value = "'-0.35%"
$ws.Range("E26").Value = "'-3.94%"
$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'0.08%"
$ws.Range("D39").Value = "'0.02834"
$ws.Range("E39").Value = "'-0.26%"
$ws.Range("D40").Value = "'0.05522"
$ws.Range("E40").Value = "'1.02%"
$ws.Range("D41").Value = "'0.007924"
$ws.Range("E41").Value = "'4.17%"
$ws.Range("E42").Value = "'1.29%"
$ws.Range("D43").Value = "'0.008985"
$ws.Range("E43").Value = "'-9.80%"
$ws.Range("D44").Value = "'0.002094"
$ws.Range("E44").Value = "'-1.78%"
$ws.Range("D45").Value = "'0.01100"
$ws.Range("E45").Value = "'-7.52%"
$ws.Range("D46").Value = "'0.00007277"
$ws.Range("E46").Value = "'8.44%"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("D48").Value = "'0.003251"
$ws.Range("E48").Value = "'8.74%"
$ws.Range("D49").Value = "'0.002282"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.08%"
